# "Fixing Report and Analysis"
# Adds two new comparison tables on the DATA sheet (below the existing
# "Number of Nodes / Genetic Algorithm(default) / Tabu(default)" block at
# C22:D29), mirroring the Brute Force/DP/GA/Tabu comparison table that
# already lives at G2:K9:
#   - G25:I32  -> Number of Nodes | Dynamic Programming | Genetic Algorithm
#   - G34:I41  -> Number of Nodes | Dynamic Programming | Tabu
# Both new tables are exact copies (values + styles) of the corresponding
# columns from the existing G2:K9 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1: Number of Nodes / Dynamic Programming / Genetic Algorithm ---
# header + 8 data rows, placed at G25:I32
$ws.Range("G2:G9").Copy($ws.Range("G25"))
$ws.Range("I2:I9").Copy($ws.Range("H25"))
$ws.Range("J2:J9").Copy($ws.Range("I25"))

# --- Table 2: Number of Nodes / Dynamic Programming / Tabu ---
# header + 8 data rows, placed at G34:I41
$ws.Range("G2:G9").Copy($ws.Range("G34"))
$ws.Range("I2:I9").Copy($ws.Range("H34"))
$ws.Range("K2:K9").Copy($ws.Range("I34"))

# Restore the selection/viewport to where the editor ended up.
$ws.Range("J29").Select()
